# ---------------------------------------------------------------------------
# Edit: Fri, Apr 24, 2020  4:06:26 PM
#
# 1) The table on slide 6 ("SOURCES OF FINANCE") switches its table style
#    from the local custom style {2936F0D0-81F0-4C4B-A4DD-D06AA56CEE3D}
#    to the built-in style {B8E27276-FEA6-4755-8D1E-4E8F46AAE9B5}.
#
# 2) The presentation's theme colours revert from the "Integral" design
#    back to the default "Office Theme" palette (the deck's two theme
#    parts end up holding each other's original colour scheme).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style -----------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{B8E27276-FEA6-4755-8D1E-4E8F46AAE9B5}")

# --- 2) Theme colours -----------------------------------------------------
function Hex2Bgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock Office Theme colours (what the deck's
# active theme reverts to).
$officeColors = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = Hex2Bgr($officeColors[$i - 1])
}
